$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("D17").Value = 44530
$ws.Range("J17").Value = 35
# Row 18
$ws.Range("D18").Value = 44365
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("P18").Value = 2000
# Row 19
$ws.Range("D19").Value = 44476
$ws.Range("J19").Value = 35
$ws.Range("M19").Value = 11429
$ws.Range("P19").Value = 1905
# Row 20
$ws.Range("D20").Value = 44411
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 11000
$ws.Range("M20").Value = 11500
$ws.Range("P20").Value = 1917
# Row 21
$ws.Range("D21").Value = 44364
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 12000
$ws.Range("P21").Value = 2000
# Row 22
$ws.Range("D22").Value = 44313
# Row 23
$ws.Range("D23").Value = 44334
$ws.Range("I23").Value = 'Segunda'
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = 9000
$ws.Range("P23").Value = 1500
# Row 24
$ws.Range("D24").Value = 44414
$ws.Range("I24").Value = 'Primera'
$ws.Range("K24").Value = 12000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 12000
$ws.Range("P24").Value = 2000
# Row 25
$ws.Range("D25").Value = 44316
$ws.Range("I25").Value = 'Segunda'
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 9000
$ws.Range("P25").Value = 1500
# Row 26
$ws.Range("D26").Value = 44466
$ws.Range("J26").Value = 30
# Row 27
$ws.Range("D27").Value = 44385
$ws.Range("J27").Value = 25
# Row 28
$ws.Range("D28").Value = 44379
$ws.Range("K28").Value = 12000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 12000
$ws.Range("P28").Value = 2000
# Row 29
$ws.Range("D29").Value = 44253
# Row 30
$ws.Range("D30").Value = 44484
$ws.Range("K30").Value = 11000
$ws.Range("L30").Value = 11000
$ws.Range("M30").Value = 11000
$ws.Range("P30").Value = 1833
# Row 31
$ws.Range("D31").Value = 44418
$ws.Range("J31").Value = 40
# Row 32
$ws.Range("D32").Value = 44427
$ws.Range("J32").Value = 20
# Row 33
$ws.Range("D33").Value = 44243
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 12000
$ws.Range("P33").Value = 2000
# Row 34
$ws.Range("D34").Value = 44280
$ws.Range("J34").Value = 30
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = 11000
$ws.Range("P34").Value = 1833
# Row 35
$ws.Range("D35").Value = 44369
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 11000
$ws.Range("M35").Value = 11500
$ws.Range("P35").Value = 1917
# Row 36
$ws.Range("D36").Value = 44397
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 12000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 12000
$ws.Range("P36").Value = 2000
# Row 37
$ws.Range("D37").Value = 44344
$ws.Range("I37").Value = 'Segunda'
$ws.Range("K37").Value = 9000
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = 9000
$ws.Range("P37").Value = 1500
# Row 38
$ws.Range("D38").Value = 44383
$ws.Range("J38").Value = 40
# Row 39
$ws.Range("D39").Value = 44267
$ws.Range("J39").Value = 50
# Row 40
$ws.Range("D40").Value = 44354
# Row 41
$ws.Range("D41").Value = 44412
$ws.Range("J41").Value = 20
# Row 42
$ws.Range("D42").Value = 44525
$ws.Range("H42").Value = 'Americana (o)'
$ws.Range("J42").Value = 35
# Row 43
$ws.Range("D43").Value = 44266
$ws.Range("H43").Value = 'Sin especificar'
$ws.Range("J43").Value = 20
# Row 44
$ws.Range("D44").Value = 44473
$ws.Range("J44").Value = 35
# Row 45
$ws.Range("D45").Value = 44433
$ws.Range("J45").Value = 20
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = 12000
$ws.Range("P45").Value = 2000
# Row 46
$ws.Range("D46").Value = 44224
$ws.Range("J46").Value = 30
$ws.Range("K46").Value = 13000
$ws.Range("L46").Value = 13000
$ws.Range("M46").Value = 13000
$ws.Range("P46").Value = 2167
# Row 47
$ws.Range("D47").Value = 44526
$ws.Range("J47").Value = 35
# Row 48
$ws.Range("D48").Value = 44186
$ws.Range("J48").Value = 8
$ws.Range("K48").Value = 11000
$ws.Range("L48").Value = 11000
$ws.Range("M48").Value = 11000
$ws.Range("P48").Value = 1833
# Row 49
$ws.Range("D49").Value = 44460
$ws.Range("J49").Value = 30
$ws.Range("K49").Value = 14000
$ws.Range("L49").Value = 14000
$ws.Range("M49").Value = 14000
$ws.Range("P49").Value = 2333
# Row 50
$ws.Range("D50").Value = 44438
$ws.Range("J50").Value = 20
# Row 51
$ws.Range("D51").Value = 44519
$ws.Range("J51").Value = 45
# Row 52
$ws.Range("D52").Value = 44392
# Row 53
$ws.Range("D53").Value = 44355
$ws.Range("H53").Value = 'Americana (o)'
$ws.Range("J53").Value = 40
# Row 54
$ws.Range("D54").Value = 44434
$ws.Range("H54").Value = 'Sin especificar'
$ws.Range("J54").Value = 25
$ws.Range("K54").Value = 12000
$ws.Range("L54").Value = 12000
$ws.Range("M54").Value = 12000
$ws.Range("P54").Value = 2000
# Row 55
$ws.Range("D55").Value = 44497
$ws.Range("J55").Value = 35
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 11000
$ws.Range("M55").Value = 10429
$ws.Range("P55").Value = 1738
# Row 56
$ws.Range("D56").Value = 44449
$ws.Range("K56").Value = 12000
$ws.Range("L56").Value = 12000
$ws.Range("M56").Value = 12000
$ws.Range("P56").Value = 2000
# Row 57
$ws.Range("D57").Value = 44358
$ws.Range("J57").Value = 40
$ws.Range("K57").Value = 11000
$ws.Range("L57").Value = 11000
$ws.Range("M57").Value = 11000
$ws.Range("P57").Value = 1833
# Row 58
$ws.Range("D58").Value = 44399
$ws.Range("J58").Value = 25
# Row 59
$ws.Range("D59").Value = 44298
$ws.Range("K59").Value = 12000
$ws.Range("L59").Value = 12000
$ws.Range("M59").Value = 12000
$ws.Range("P59").Value = 2000
# Row 60
$ws.Range("D60").Value = 44482
$ws.Range("J60").Value = 20
$ws.Range("K60").Value = 11000
$ws.Range("L60").Value = 11000
$ws.Range("M60").Value = 11000
$ws.Range("P60").Value = 1833
# Row 61
$ws.Range("D61").Value = 44250
$ws.Range("J61").Value = 40
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = 10000
$ws.Range("P61").Value = 1667
# Row 62
$ws.Range("D62").Value = 44273
$ws.Range("J62").Value = 20
# Row 63
$ws.Range("D63").Value = 44386
$ws.Range("J63").Value = 50
# Row 64
$ws.Range("H64").Value = 'Americana (o)'
$ws.Range("J64").Value = 65
# Row 65
$ws.Range("D65").Value = 44435
$ws.Range("H65").Value = 'Sin especificar'
$ws.Range("J65").Value = 70
# Row 66
$ws.Range("D66").Value = 44277
$ws.Range("J66").Value = 20
# Row 67
$ws.Range("D67").Value = 44442
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = 12000
$ws.Range("P67").Value = 2000
# Row 68
$ws.Range("D68").Value = 44516
$ws.Range("J68").Value = 40
$ws.Range("K68").Value = 11000
$ws.Range("L68").Value = 11000
$ws.Range("M68").Value = 11000
$ws.Range("P68").Value = 1833
# Row 69
$ws.Range("D69").Value = 44175
$ws.Range("J69").Value = 20
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = 12000
$ws.Range("P69").Value = 2000
# Row 70
$ws.Range("D70").Value = 44168
$ws.Range("J70").Value = 15
$ws.Range("K70").Value = 13000
$ws.Range("L70").Value = 13000
$ws.Range("M70").Value = 13000
$ws.Range("P70").Value = 2167
# Row 71
$ws.Range("D71").Value = 44203
$ws.Range("K71").Value = 14000
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = 14000
$ws.Range("P71").Value = 2333
# Row 72
$ws.Range("D72").Value = 44475
$ws.Range("J72").Value = 20
$ws.Range("K72").Value = 12000
$ws.Range("L72").Value = 12000
$ws.Range("M72").Value = 12000
$ws.Range("P72").Value = 2000
# Row 73
$ws.Range("D73").Value = 44483
$ws.Range("J73").Value = 25
$ws.Range("K73").Value = 11000
$ws.Range("L73").Value = 11000
$ws.Range("M73").Value = 11000
$ws.Range("P73").Value = 1833
# Row 74
$ws.Range("D74").Value = 44217
$ws.Range("J74").Value = 20
$ws.Range("K74").Value = 12000
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = 12000
$ws.Range("P74").Value = 2000
# Row 75
$ws.Range("D75").Value = 44162
$ws.Range("J75").Value = 40
$ws.Range("K75").Value = 13000
$ws.Range("L75").Value = 13000
$ws.Range("M75").Value = 13000
$ws.Range("P75").Value = 2167
# Row 76
$ws.Range("D76").Value = 44357
$ws.Range("J76").Value = 25
$ws.Range("K76").Value = 11500
$ws.Range("L76").Value = 11500
$ws.Range("M76").Value = 11500
$ws.Range("P76").Value = 1917
# Row 77
$ws.Range("D77").Value = 44333
$ws.Range("J77").Value = 15
# Row 78
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 20
$ws.Range("K78").Value = 11000
$ws.Range("L78").Value = 11000
$ws.Range("M78").Value = 11000
$ws.Range("P78").Value = 1833
# Row 79
$ws.Range("D79").Value = 44320
$ws.Range("I79").Value = 'Segunda'
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = 9000
$ws.Range("P79").Value = 1500
# Row 80
$ws.Range("D80").Value = 44252
$ws.Range("J80").Value = 30
$ws.Range("K80").Value = 11000
$ws.Range("L80").Value = 11000
$ws.Range("M80").Value = 11000
$ws.Range("P80").Value = 1833
# Row 81
$ws.Range("D81").Value = 44467
$ws.Range("K81").Value = 12000
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = 12000
$ws.Range("P81").Value = 2000
# Row 82
$ws.Range("D82").Value = 44264
$ws.Range("J82").Value = 40
$ws.Range("K82").Value = 11000
$ws.Range("L82").Value = 11000
$ws.Range("M82").Value = 11000
$ws.Range("P82").Value = 1833
# Row 83
$ws.Range("D83").Value = 44167
$ws.Range("J83").Value = 15
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = 13000
$ws.Range("P83").Value = 2167
# Row 84
$ws.Range("D84").Value = 44390
$ws.Range("J84").Value = 40
# Row 85
$ws.Range("D85").Value = 44291
$ws.Range("I85").Value = 'Primera'
$ws.Range("J85").Value = 20
$ws.Range("K85").Value = 11000
$ws.Range("L85").Value = 11000
$ws.Range("M85").Value = 11000
$ws.Range("P85").Value = 1833
# Row 86
$ws.Range("D86").Value = 44326
$ws.Range("I86").Value = 'Segunda'
$ws.Range("J86").Value = 25
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = 9000
$ws.Range("P86").Value = 1500
# Row 87
$ws.Range("D87").Value = 44302
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 12000
$ws.Range("M87").Value = 12000
$ws.Range("P87").Value = 2000
# Row 88
$ws.Range("D88").Value = 44292
$ws.Range("J88").Value = 40
$ws.Range("K88").Value = 11000
$ws.Range("L88").Value = 11000
$ws.Range("M88").Value = 11000
$ws.Range("P88").Value = 1833
# Row 89
$ws.Range("D89").Value = 44308
$ws.Range("J89").Value = 15
$ws.Range("K89").Value = 12000
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = 12000
$ws.Range("P89").Value = 2000
# Row 90
$ws.Range("D90").Value = 44498
$ws.Range("I90").Value = 'Primera'
$ws.Range("J90").Value = 40
$ws.Range("K90").Value = 10000
$ws.Range("L90").Value = 11000
$ws.Range("M90").Value = 10500
$ws.Range("P90").Value = 1750
# Row 91
$ws.Range("D91").Value = 44321
$ws.Range("I91").Value = 'Segunda'
$ws.Range("J91").Value = 15
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 9000
$ws.Range("M91").Value = 9000
$ws.Range("P91").Value = 1500
# Row 92
$ws.Range("D92").Value = 44477
$ws.Range("J92").Value = 35
$ws.Range("K92").Value = 11000
$ws.Range("L92").Value = 11000
$ws.Range("M92").Value = 11000
$ws.Range("P92").Value = 1833
# Row 93
$ws.Range("D93").Value = 44487
$ws.Range("J93").Value = 30
# Row 94
$ws.Range("D94").Value = 44452
$ws.Range("J94").Value = 25
$ws.Range("K94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = 12000
$ws.Range("P94").Value = 2000
# Row 95
$ws.Range("D95").Value = 44505
$ws.Range("J95").Value = 50
$ws.Range("K95").Value = 9000
$ws.Range("L95").Value = 10000
$ws.Range("M95").Value = 9500
$ws.Range("P95").Value = 1583
# Row 96
$ws.Range("D96").Value = 44204
$ws.Range("K96").Value = 14000
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = 14500
$ws.Range("P96").Value = 2417
# Row 97
$ws.Range("D97").Value = 44260
$ws.Range("K97").Value = 11000
$ws.Range("L97").Value = 11000
$ws.Range("M97").Value = 11000
$ws.Range("P97").Value = 1833
# Row 98
$ws.Range("D98").Value = 44306
$ws.Range("J98").Value = 40
$ws.Range("K98").Value = 12000
$ws.Range("L98").Value = 12000
$ws.Range("M98").Value = 12000
$ws.Range("P98").Value = 2000
# Row 99
$ws.Range("D99").Value = 44509
$ws.Range("J99").Value = 20
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = 10000
$ws.Range("P99").Value = 1667
# Row 100
$ws.Range("D100").Value = 44278
$ws.Range("J100").Value = 40
# Row 101
$ws.Range("D101").Value = 44494
# Row 102
$ws.Range("D102").Value = 44469
# Row 103
$ws.Range("D103").Value = 44518
$ws.Range("J103").Value = 25
# Row 104
$ws.Range("D104").Value = 44446
$ws.Range("J104").Value = 50
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 12000
$ws.Range("M104").Value = 12000
$ws.Range("P104").Value = 2000
# Row 105
$ws.Range("D105").Value = 44463
$ws.Range("J105").Value = 40
$ws.Range("K105").Value = 14000
$ws.Range("L105").Value = 14000
$ws.Range("M105").Value = 14000
$ws.Range("P105").Value = 2333
# Row 106
$ws.Range("D106").Value = 44245
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 20
$ws.Range("K106").Value = 10000
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = 10000
$ws.Range("P106").Value = 1667
# Row 107
$ws.Range("D107").Value = 44323
$ws.Range("I107").Value = 'Segunda'
$ws.Range("K107").Value = 9000
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = 9000
$ws.Range("P107").Value = 1500
# Row 108
$ws.Range("D108").Value = 44481
$ws.Range("J108").Value = 40
$ws.Range("K108").Value = 11000
$ws.Range("L108").Value = 11000
$ws.Range("M108").Value = 11000
$ws.Range("P108").Value = 1833
# Row 109
$ws.Range("D109").Value = 44229
$ws.Range("J109").Value = 50
$ws.Range("K109").Value = 13000
$ws.Range("L109").Value = 13000
$ws.Range("M109").Value = 13000
$ws.Range("P109").Value = 2167
# Row 110
$ws.Range("D110").Value = 44417
$ws.Range("J110").Value = 20
$ws.Range("K110").Value = 12000
$ws.Range("L110").Value = 12000
$ws.Range("M110").Value = 12000
$ws.Range("P110").Value = 2000
# Row 111
$ws.Range("D111").Value = 44523
$ws.Range("J111").Value = 40
$ws.Range("K111").Value = 11000
$ws.Range("L111").Value = 11000
$ws.Range("M111").Value = 11000
$ws.Range("P111").Value = 1833
# Row 112
$ws.Range("D112").Value = 44462
$ws.Range("J112").Value = 30
$ws.Range("K112").Value = 12000
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = 12000
$ws.Range("P112").Value = 2000
# Row 113
$ws.Range("D113").Value = 44259
$ws.Range("J113").Value = 20
$ws.Range("K113").Value = 11000
$ws.Range("L113").Value = 11000
$ws.Range("M113").Value = 11000
$ws.Range("P113").Value = 1833
# Row 114
$ws.Range("D114").Value = 44406
$ws.Range("J114").Value = 30
# Row 115
$ws.Range("D115").Value = 44295
$ws.Range("J115").Value = 40
# Row 116
$ws.Range("D116").Value = 44270
$ws.Range("J116").Value = 20
$ws.Range("K116").Value = 12000
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = 12000
$ws.Range("P116").Value = 2000
# Row 117
$ws.Range("I117").Value = 'Primera'
$ws.Range("J117").Value = 30
$ws.Range("K117").Value = 11000
$ws.Range("L117").Value = 11000
$ws.Range("M117").Value = 11000
$ws.Range("P117").Value = 1833
# Row 118
$ws.Range("D118").Value = 44299
$ws.Range("I118").Value = 'Segunda'
$ws.Range("J118").Value = 20
$ws.Range("K118").Value = 9000
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 9000
$ws.Range("P118").Value = 1500
# Row 119
$ws.Range("D119").Value = 44257
$ws.Range("J119").Value = 50
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 10000
$ws.Range("M119").Value = 10000
$ws.Range("P119").Value = 1667
# Row 120
$ws.Range("D120").Value = 44372
$ws.Range("J120").Value = 40
$ws.Range("K120").Value = 11000
$ws.Range("M120").Value = 11500
$ws.Range("P120").Value = 1917
# Row 121
$ws.Range("D121").Value = 44403
$ws.Range("J121").Value = 20
# Row 122
$ws.Range("D122").Value = 44169
# Row 123
$ws.Range("D123").Value = 44376
$ws.Range("J123").Value = 40
$ws.Range("K123").Value = 12000
$ws.Range("M123").Value = 12000
$ws.Range("P123").Value = 2000
# Row 124
$ws.Range("D124").Value = 44474
$ws.Range("J124").Value = 50
$ws.Range("K124").Value = 11000
$ws.Range("M124").Value = 11500
$ws.Range("P124").Value = 1917
# Row 125
$ws.Range("D125").Value = 44421
$ws.Range("J125").Value = 40
# Row 126
$ws.Range("D126").Value = 44242
$ws.Range("J126").Value = 30
# Row 127
$ws.Range("D127").Value = 44239
$ws.Range("J127").Value = 50
# Row 128
$ws.Range("D128").Value = 44448
$ws.Range("J128").Value = 25
$ws.Range("K128").Value = 12000
$ws.Range("L128").Value = 12000
$ws.Range("M128").Value = 12000
$ws.Range("P128").Value = 2000
# Row 129
$ws.Range("D129").Value = 44362
$ws.Range("I129").Value = 'Primera'
$ws.Range("J129").Value = 40
# Row 130
$ws.Range("D130").Value = 44210
$ws.Range("I130").Value = 'Segunda'
$ws.Range("J130").Value = 20
$ws.Range("K130").Value = 11000
$ws.Range("L130").Value = 11000
$ws.Range("M130").Value = 11000
$ws.Range("P130").Value = 1833
# Row 131
$ws.Range("D131").Value = 44176
$ws.Range("J131").Value = 40
# Row 132
$ws.Range("D132").Value = 44301
$ws.Range("J132").Value = 25
# Row 133
$ws.Range("D133").Value = 44407
$ws.Range("J133").Value = 40
$ws.Range("K133").Value = 12000
$ws.Range("L133").Value = 12000
$ws.Range("M133").Value = 12000
$ws.Range("P133").Value = 2000
# Row 134
$ws.Range("D134").Value = 44284
$ws.Range("J134").Value = 15
$ws.Range("K134").Value = 11000
$ws.Range("L134").Value = 11000
$ws.Range("M134").Value = 11000
$ws.Range("P134").Value = 1833
# Row 135
$ws.Range("I135").Value = 'Primera'
$ws.Range("J135").Value = 20
$ws.Range("K135").Value = 12000
$ws.Range("L135").Value = 12000
$ws.Range("M135").Value = 12000
$ws.Range("P135").Value = 2000
# Row 136
$ws.Range("D136").Value = 44341
$ws.Range("I136").Value = 'Segunda'
$ws.Range("J136").Value = 30
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 9000
$ws.Range("P136").Value = 1500
# Row 137
$ws.Range("D137").Value = 44504
$ws.Range("J137").Value = 25
$ws.Range("K137").Value = 11000
$ws.Range("L137").Value = 11000
$ws.Range("M137").Value = 11000
$ws.Range("P137").Value = 1833
# Row 138
$ws.Range("D138").Value = 44350
$ws.Range("J138").Value = 20
$ws.Range("K138").Value = 12000
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 12000
$ws.Range("P138").Value = 2000
# Row 139
$ws.Range("D139").Value = 44329
$ws.Range("J139").Value = 30
# Row 140
$ws.Range("D140").Value = 44522
$ws.Range("J140").Value = 25
$ws.Range("K140").Value = 11000
$ws.Range("L140").Value = 11000
$ws.Range("M140").Value = 11000
$ws.Range("P140").Value = 1833
# Row 141
$ws.Range("D141").Value = 44246
$ws.Range("J141").Value = 30
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 10000
$ws.Range("M141").Value = 10000
$ws.Range("P141").Value = 1667
# Row 142
$ws.Range("I142").Value = 'Primera'
$ws.Range("K142").Value = 11000
$ws.Range("L142").Value = 11000
$ws.Range("M142").Value = 11000
$ws.Range("P142").Value = 1833
# Row 143
$ws.Range("D143").Value = 44491
$ws.Range("I143").Value = 'Segunda'
$ws.Range("J143").Value = 25
$ws.Range("K143").Value = 10000
$ws.Range("L143").Value = 10000
$ws.Range("M143").Value = 10000
$ws.Range("P143").Value = 1667
# Row 144
$ws.Range("D144").Value = 44166
$ws.Range("K144").Value = 13000
$ws.Range("M144").Value = 13000
$ws.Range("P144").Value = 2167
# Row 145
$ws.Range("D145").Value = 44225
$ws.Range("J145").Value = 40
$ws.Range("L145").Value = 13000
$ws.Range("M145").Value = 12500
$ws.Range("P145").Value = 2083
# Row 146
$ws.Range("D146").Value = 44447
$ws.Range("J146").Value = 20
# Row 147
$ws.Range("D147").Value = 44425
$ws.Range("J147").Value = 50
# Row 148
$ws.Range("D148").Value = 44348
$ws.Range("J148").Value = 40
$ws.Range("K148").Value = 12000
$ws.Range("L148").Value = 12000
$ws.Range("M148").Value = 12000
$ws.Range("P148").Value = 2000
# Row 149
$ws.Range("D149").Value = 44322
$ws.Range("J149").Value = 20
$ws.Range("K149").Value = 11000
$ws.Range("L149").Value = 11000
$ws.Range("M149").Value = 11000
$ws.Range("P149").Value = 1833
# Row 150
$ws.Range("D150").Value = 44495
$ws.Range("K150").Value = 10000
$ws.Range("L150").Value = 10000
$ws.Range("M150").Value = 10000
$ws.Range("P150").Value = 1667
# Row 151
$ws.Range("D151").Value = 44232
$ws.Range("J151").Value = 40
$ws.Range("K151").Value = 12000
$ws.Range("L151").Value = 12000
$ws.Range("M151").Value = 12000
$ws.Range("P151").Value = 2000
# Row 152
$ws.Range("I152").Value = 'Primera'
$ws.Range("J152").Value = 30
$ws.Range("K152").Value = 11000
$ws.Range("L152").Value = 11000
$ws.Range("M152").Value = 11000
$ws.Range("P152").Value = 1833
# Row 153
$ws.Range("D153").Value = 44327
$ws.Range("I153").Value = 'Segunda'
$ws.Range("K153").Value = 9000
$ws.Range("L153").Value = 9000
$ws.Range("M153").Value = 9000
$ws.Range("P153").Value = 1500
# Row 154
$ws.Range("D154").Value = 44161
$ws.Range("J154").Value = 20
$ws.Range("K154").Value = 14000
$ws.Range("L154").Value = 14000
$ws.Range("M154").Value = 14000
$ws.Range("P154").Value = 2333
# Row 155
$ws.Range("D155").Value = 44468
$ws.Range("J155").Value = 15
# Row 156
$ws.Range("D156").Value = 44238
$ws.Range("J156").Value = 30
# Row 157
$ws.Range("D157").Value = 44236
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 40
$ws.Range("K157").Value = 12000
$ws.Range("L157").Value = 12000
$ws.Range("M157").Value = 12000
$ws.Range("P157").Value = 2000
# Row 158
$ws.Range("D158").Value = 44340
$ws.Range("I158").Value = 'Segunda'
$ws.Range("J158").Value = 15
$ws.Range("K158").Value = 9000
$ws.Range("L158").Value = 9000
$ws.Range("M158").Value = 9000
$ws.Range("P158").Value = 1500
# Row 159
$ws.Range("I159").Value = 'Primera'
$ws.Range("J159").Value = 30
$ws.Range("K159").Value = 11000
$ws.Range("L159").Value = 11000
$ws.Range("M159").Value = 11000
$ws.Range("P159").Value = 1833
# Row 160
$ws.Range("D160").Value = 44330
$ws.Range("I160").Value = 'Segunda'
$ws.Range("J160").Value = 20
$ws.Range("K160").Value = 9000
$ws.Range("L160").Value = 9000
$ws.Range("M160").Value = 9000
$ws.Range("P160").Value = 1500
# Row 161
$ws.Range("D161").Value = 44432
$ws.Range("J161").Value = 45
$ws.Range("K161").Value = 12000
$ws.Range("M161").Value = 12000
$ws.Range("P161").Value = 2000
# Row 162
$ws.Range("D162").Value = 44181
$ws.Range("J162").Value = 20
$ws.Range("K162").Value = 11000
$ws.Range("M162").Value = 11500
$ws.Range("P162").Value = 1917
# Row 163
$ws.Range("D163").Value = 44271
$ws.Range("J163").Value = 40
# Row 164
$ws.Range("D164").Value = 44400
$ws.Range("J164").Value = 50
$ws.Range("K164").Value = 12000
$ws.Range("M164").Value = 12000
$ws.Range("P164").Value = 2000
# Row 165
$ws.Range("D165").Value = 44309
$ws.Range("J165").Value = 40
$ws.Range("K165").Value = 11000
$ws.Range("L165").Value = 12000
$ws.Range("M165").Value = 11500
$ws.Range("P165").Value = 1917
# Row 166
$ws.Range("D166").Value = 44508
$ws.Range("J166").Value = 25

# Row 167 (new)
$ws.Range("A167").Value = 4
$ws.Range("B167").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C167").Value = 'Los Lagos'
$ws.Range("D167").Value = 44201
$ws.Range("E167").Value = 10
$ws.Range("F167").Value = 100112017
$ws.Range("G167").Value = 'Apio'
$ws.Range("H167").Value = 'Americana (o)'
$ws.Range("I167").Value = 'Primera'
$ws.Range("J167").Value = 30
$ws.Range("K167").Value = 10000
$ws.Range("L167").Value = 10000
$ws.Range("M167").Value = 10000
$ws.Range("N167").Value = '$/docena de matas'
$ws.Range("O167").Value = 'Región de Coquimbo'
$ws.Range("P167").Value = 1667
$ws.Range("Q167").Value = 6
$ws.Range("R167").Value = 'Hortaliza'

# Match date number format/style for new row
$ws.Range("D167").NumberFormat = $ws.Range("D166").NumberFormat
